$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: was TC_01 ("no parameter" / Fail) -> now TC_02 ("Geopolitical ID" / Pass) ---
$ws.Range("A2").Value = "TC_02"
$ws.Range("B2").Value = "Verify that the StProvStd service is successfully retrieving the records when passing the Geopolitical ID in URI"
$ws.Range("D2").Value = "NA"
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = "Pass"
$ws.Range("K2").Value = "Total number of records matching between DB & Response: 0, below are the test steps for this test case"
# Result flipped from Fail to Pass -> font turns from red to green, matching rows 3 & 4
# (VBA ColorIndex 10 = RGB 008000 dark-green, i.e. xlsx raw indexed=17)
$ws.Range("A2:K2").Font.ColorIndex = 10

# --- Row 3: was TC_02 ("Geopolitical ID") -> now TC_03 ("Geopolitical ID, targetDate and endDate") ---
$ws.Range("A3").Value = "TC_03"
$ws.Range("B3").Value = "Verify that the StProvStd service is successfully retrieving the records when passing the Geopolitical ID, targetDate and endDate in URI"
$ws.Range("K3").Value = "Total number of records matching between DB & Response: 0, below are the test steps for this test case"

# --- Row 4: was blank TestCase ID ("Geopolitical ID" dup row) -> now TC_12 (invalid stProvCd) ---
$ws.Range("A4").Value = "TC_12"
$ws.Range("B4").Value = "Verify that the StProvStd service is not retrieving the records when passing the invalid stProvCd and valid orgStdCd  in URI"
$ws.Range("G4").Value = "Success"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "200"
$ws.Range("I4").Value = ""
$ws.Range("K4").Value = "Total number of records matching between DB & Response: 0, below are the test steps for this test case"
